# "fix the bug of index" - correct the shifted Count values in the second
# (Adapted/Original) table on Sheet1: what had been recorded in C56 really
# belonged one row down (C57 is the SUM row), so C55 picks up the missing
# value (48) and C56 is corrected to its proper total (72).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C55").Value = 48
$ws.Range("C56").Value = 72

# Recalculate so the SUM(C50:C56) row and any dependent chart series pick up
# the corrected values.
$excel.CalculateFull()

# Reflect where the user ended up looking after making the fix: scrolled
# down to row 45 and with C57 (the corrected total) selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 45
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C57").Select()
